$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking
# strings (e.g. "32.65") are stored as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '34.023.09'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '1.780.03'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '227.13'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("E6").Value = '  -0.94%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '32.65'
$ws.Range("E8").Value = '  +2.39%  '
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("D10").Value = '0.0712'
$ws.Range("E10").Value = '  -2.34%  '
$ws.Range("D11").Value = '0.0936'
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").Value = '2.035.67'
$ws.Range("E12").Value = '  -1.70%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.790.13'
$ws.Range("E13").Value = '  -1.22%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '10.93'
$ws.Range("E14").Value = '  +0.88%  '
$ws.Range("D15").Value = '34.006.98'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").Value = '0.619'
$ws.Range("E16").Value = '  -3.10%  '
$ws.Range("E17").Value = '  -3.95%  '
$ws.Range("D18").Value = '67.75'
$ws.Range("E18").Value = '  -2.12%  '
$ws.Range("D19").Value = '244.54'
$ws.Range("E19").Value = '  -1.90%  '
$ws.Range("D20").Value = '0.0₃0784'
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '10.71'
$ws.Range("E22").Value = '  -2.27%  '
$ws.Range("D23").Value = '4.08'
$ws.Range("E23").Value = '  -3.73%  '
$ws.Range("E24").Value = '  -3.22%  '
$ws.Range("D25").Value = '160.12'
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = '16.30'
$ws.Range("E26").Value = '  -1.61%  '
$ws.Range("D27").Value = '7.04'
$ws.Range("E27").Value = '  -2.02%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").Value = '1.23'
$ws.Range("E30").Value = '  +1.77%  '
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("D32").Value = '3.63'
$ws.Range("E32").Value = '  -3.30%  '
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("D34").Value = '1.81'
$ws.Range("E34").Value = '  -3.59%  '
$ws.Range("D35").Value = '1.388.00'
$ws.Range("E35").Value = '  -2.94%  '
$ws.Range("D36").Value = '0.649'
$ws.Range("E36").Value = '  +2.16%  '
$ws.Range("E37").Value = '  -1.43%  '
$ws.Range("E38").Value = '  -1.35%  '
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("D40").Value = '2.19'
$ws.Range("E40").Value = '  +1.88%  '
$ws.Range("D41").Value = '0.913'
$ws.Range("E41").Value = '  -4.40%  '
$ws.Range("E42").Value = '  -2.43%  '
$ws.Range("D43").Value = '77.75'
$ws.Range("E43").Value = '  -3.76%  '
$ws.Range("D44").Value = '0.0₆0139'
$ws.Range("E44").Value = '  +14.45%  '
$ws.Range("D46").Value = '12.72'
$ws.Range("E46").Value = '  +8.09%  '
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D48").Value = '107.53'
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("D50").Value = '1.936.29'
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("E51").Value = '  +0.11%  '

# Restore default styling (no explicit style was present on these cells originally)
$ws.Range("D2:E51").Style = "Normal"
